# Add actual work-hours for the features assigned to Jabesi (sprint backlog update).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (login/logout task): record actual hours spent in Week 1 / Week 2,
# and bump the initial estimate from 6.5 to 8.
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 3

# Row 6 (admin - add nurses task): assign Rahman.
$ws.Range("D6").Value = "Rahman"

# Rows 8 & 9 (design tasks): assign the team and mark estimate as not applicable.
$ws.Range("D8").Value = "Jabesi/Ahmad/Rahman"
$ws.Range("E8").Value = "N/A"
$ws.Range("D9").Value = "Jabesi/Ahmad/Rahman"
$ws.Range("E9").Value = "N/A"

# Move the active selection to B7.
$ws.Range("B7").Select()

# Resize the burndown chart (its bottom edge was dragged further down).
$co = $ws.ChartObjects(1)
$co.Height = 187.5
